$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# productname: append "-1st" suffix to the product code on both sheets
$wsInput.Range("B1").Value = "2615-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-MIS-1st"
$wsOutput.Range("B1").Value = "2615-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-MIS-1st"

# shortname: change from numeric 2615 to text "261e"
$wsInput.Range("B2").Value = "261e"

# Move the active selection from B15 to B7 on the input sheet
$wsInput.Range("B7").Select()
